$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Workbook window geometry (best effort - mirrors the diff's bookViews
#    windowHeight/windowWidth/xWindow change)
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Height = 792
$win.Width = 1452
$win.Left = 1434
$win.Top = -6

# ---------------------------------------------------------------------------
# 2) Header / year labels: 2006 -> 2012
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------------
# 3) Updated factor data (Average Values column E, Ridership Effect column H)
#    Column F (Year2 average values) is unchanged in the diff.
# ---------------------------------------------------------------------------
$ws1.Range("E8").Value = 438984
$ws1.Range("H8").Value = 39722.20194999999

$ws1.Range("E9").Value = 0.610122203
$ws1.Range("H9").Value = 10312.505709

$ws1.Range("E10").Value = 136089.92
$ws1.Range("H10").Value = 18014.87546

$ws1.Range("E11").Value = 3.032949666
$ws1.Range("H11").Value = -4395.5385829

$ws1.Range("E12").Value = 3.9458
$ws1.Range("H12").Value = -38125.99518250001

$ws1.Range("E13").Value = 22730.86
$ws1.Range("H13").Value = -15907.302115

$ws1.Range("E14").Value = 8.49
$ws1.Range("H14").Value = 15063.752053

$ws1.Range("E15").Value = 2.3
$ws1.Range("H15").Value = 3965.9811241

# Row 16 (Years Since Ride-hail Start): E16/F16 stay blank, H16 unchanged value.
# Row 17 (Bike Share) and 18 (Electric Scooters): values unchanged (0/0/0).

# New Reporters row 19: H19 gains an explicit 0 value (previously blank).
$ws1.Range("H19").Value = 0

# Totals
$ws1.Range("E20").Value = 990271.8005
$ws1.Range("E21").Value = 1029272

# ---------------------------------------------------------------------------
# 4) Formulas: percentage-difference formulas drop the "*100" multiplier
#    (the cells are reformatted as true percentages instead of raw numbers).
# ---------------------------------------------------------------------------
for ($r = 8; $r -le 21; $r++) {
    $ws1.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
for ($r = 8; $r -le 19; $r++) {
    $ws1.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}
# I20 and I21 keep referencing the G column directly (formula text unchanged).

# ---------------------------------------------------------------------------
# 5) Number formatting: E/F/H columns -> #,##0.00 ; G/I columns -> 0.00%
#    (borders/fonts/alignment stay exactly as they were - only the number
#    format id changes, matching the underlying style diff.)
# ---------------------------------------------------------------------------
$ws1.Range("E8:F18").NumberFormat = "#,##0.00"
$ws1.Range("H8:H18").NumberFormat = "#,##0.00"
$ws1.Range("G8:G18").NumberFormat = "0.00%"
$ws1.Range("I8:I18").NumberFormat = "0.00%"

$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"
$ws1.Range("G19").NumberFormat = "0.00%"
$ws1.Range("I19").NumberFormat = "0.00%"

$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"
$ws1.Range("G20").NumberFormat = "0.00%"
$ws1.Range("I20").NumberFormat = "0.00%"

$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("G21").NumberFormat = "0.00%"
$ws1.Range("I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 6) Sheet view: drop the frozen topLeftCell scroll position and move the
#    active selection from K20 to H21.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Application.ActiveWindow().ScrollRow = 1
$ws1.Application.ActiveWindow().ScrollColumn = 1
$ws1.Range("H21").Select() | Out-Null
